$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing rows 2-11: remove "images\" segment from the image path in column D
for ($r = 2; $r -le 11; $r++) {
    $cell = $ws.Cells.Item($r, 4)
    $current = $cell.Value()
    $updated = $current -replace "images\\", ""
    $cell.Value = $updated
}

# Append new rows 12-15
$newRows = @(
    @{
        A = "Gianpiero Lambiase Is the Man in Max Verstappen’s Ear"
        B = "Nov. 23, 2023"
        C = "The race engineer gives advice to the driver on the radio during races. Their relationship has been described as a long marriage."
        D = "output\25sp-dhabi-lambiase-inyt-02-hvqf-threeByTwoSmallAt2X.png"
        E = 0
        F = $false
    },
    @{
        A = "On Second Thought, Haas Is Staying Put"
        B = "Nov. 23, 2023"
        C = "Gene Haas once considered selling the team, but changes in F1 have convinced him not to."
        D = "output\25sp-dhabi-haas-inyt-01-wphb-threeByTwoSmallAt2X.png"
        E = 0
        F = $false
    },
    @{
        A = "McLaren Is No Longer Caught in the Middle of Formula 1"
        B = "Nov. 23, 2023"
        C = "The team started the year as another lackluster midfield team, but big changes have moved it to No. 4."
        D = "output\25sp-dhabi-mclaren-inyt-01-hlcq-threeByTwoSmallAt2X.png"
        E = 1
        F = $false
    },
    @{
        A = "The Power and Speed of Jannik Sinner"
        B = "Nov. 10, 2023"
        C = "Sinner, 22, has dominating talent and has already beaten many of tennis’s top players."
        D = "output\11sp-atp-sinner-inyt3-glvf-threeByTwoSmallAt2X.png"
        E = 0
        F = $false
    }
)

$row = 12
foreach ($item in $newRows) {
    $ws.Cells.Item($row, 1).Value = $item.A
    $ws.Cells.Item($row, 2).Value = $item.B
    $ws.Cells.Item($row, 3).Value = $item.C
    $ws.Cells.Item($row, 4).Value = $item.D
    $ws.Cells.Item($row, 5).Value = $item.E
    $ws.Cells.Item($row, 6).Value = $item.F
    $row++
}
